$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2307692307692308
$ws.Range("C2").Value = 0.4743589743589743
$ws.Range("J2").Value = 0.01602564102564102
$ws.Range("P2").Value = 0.1794871794871795
$ws.Range("S2").Value = 0.09935897435897435
# Row 3
$ws.Range("J3").Value = 0.07333333333333333
$ws.Range("P3").Value = 0.76
$ws.Range("S3").Value = 0.1666666666666667
# Row 4
$ws.Range("J4").Value = 0.05
$ws.Range("P4").Value = 0.65
$ws.Range("S4").Value = 0.3
# Row 6
$ws.Range("B6").Value = 0.06976744186046512
$ws.Range("F6").Value = 0.05813953488372093
$ws.Range("J6").Value = 0.3313953488372093
$ws.Range("O6").Value = 0.01744186046511628
$ws.Range("Q6").Value = 0.09883720930232558
$ws.Range("R6").Value = 0.04069767441860465
$ws.Range("S6").Value = 0.3837209302325582
# Row 7
$ws.Range("B7").Value = 0.1428571428571428
$ws.Range("D7").Value = 0.02597402597402598
$ws.Range("F7").Value = 0.03246753246753246
$ws.Range("J7").Value = 0.1948051948051948
$ws.Range("O7").Value = 0.03246753246753246
$ws.Range("Q7").Value = 0.1038961038961039
$ws.Range("R7").Value = 0.05844155844155844
$ws.Range("S7").Value = 0.4090909090909091
# Row 8
$ws.Range("B8").Value = 0.1187648456057007
$ws.Range("D8").Value = 0.02137767220902613
$ws.Range("F8").Value = 0.04513064133016627
$ws.Range("J8").Value = 0.1401425178147268
$ws.Range("O8").Value = 0.01900237529691211
$ws.Range("Q8").Value = 0.1900237529691211
$ws.Range("R8").Value = 0.09026128266033254
$ws.Range("S8").Value = 0.3752969121140142
# Row 9
$ws.Range("B9").Value = 0.1284403669724771
$ws.Range("D9").Value = 0.03669724770642202
$ws.Range("F9").Value = 0.08256880733944955
$ws.Range("J9").Value = 0.09174311926605505
$ws.Range("O9").Value = 0.04587155963302753
$ws.Range("Q9").Value = 0.1743119266055046
$ws.Range("R9").Value = 0.02752293577981652
$ws.Range("S9").Value = 0.4128440366972477
# Row 10
$ws.Range("B10").Value = 0.1198237885462555
$ws.Range("D10").Value = 0.02114537444933921
$ws.Range("F10").Value = 0.05991189427312775
$ws.Range("J10").Value = 0.1365638766519824
$ws.Range("O10").Value = 0.01938325991189427
$ws.Range("Q10").Value = 0.2537444933920705
$ws.Range("R10").Value = 0.07488986784140969
$ws.Range("S10").Value = 0.3145374449339207
# Row 11
$ws.Range("G11").Value = 0.2027649769585254
$ws.Range("J11").Value = 0.07834101382488479
$ws.Range("K11").Value = 0.1935483870967742
$ws.Range("L11").Value = 0.5253456221198156
# Row 12
$ws.Range("G12").Value = 0.7622950819672131
$ws.Range("J12").Value = 0.1885245901639344
$ws.Range("L12").Value = 0.03278688524590164
$ws.Range("S12").Value = 0.01639344262295082
# Row 13
$ws.Range("G13").Value = 0.6785714285714286
$ws.Range("J13").Value = 0.2857142857142857
$ws.Range("S13").Value = 0.03571428571428571
# Row 15
$ws.Range("F15").Value = 0.02298850574712644
$ws.Range("H15").Value = 0.1609195402298851
$ws.Range("I15").Value = 0.05747126436781609
$ws.Range("J15").Value = 0.4022988505747127
$ws.Range("K15").Value = 0.04597701149425287
$ws.Range("O15").Value = 0.03448275862068965
$ws.Range("S15").Value = 0.2758620689655172
# Row 16
$ws.Range("F16").Value = 0.01530612244897959
$ws.Range("H16").Value = 0.2040816326530612
$ws.Range("I16").Value = 0.06122448979591837
$ws.Range("J16").Value = 0.4387755102040816
$ws.Range("K16").Value = 0.09183673469387756
$ws.Range("M16").Value = 0.00510204081632653
$ws.Range("O16").Value = 0.04591836734693878
$ws.Range("S16").Value = 0.1377551020408163
# Row 17
$ws.Range("F17").Value = 0.02184466019417476
$ws.Range("H17").Value = 0.1868932038834951
$ws.Range("I17").Value = 0.05825242718446602
$ws.Range("J17").Value = 0.4951456310679612
$ws.Range("K17").Value = 0.08009708737864078
$ws.Range("M17").Value = 0.004854368932038835
$ws.Range("O17").Value = 0.0412621359223301
$ws.Range("S17").Value = 0.1116504854368932
# Row 18
$ws.Range("F18").Value = 0.01428571428571429
$ws.Range("H18").Value = 0.1642857142857143
$ws.Range("I18").Value = 0.07142857142857142
$ws.Range("J18").Value = 0.5
$ws.Range("K18").Value = 0.06428571428571428
$ws.Range("O18").Value = 0.05714285714285714
$ws.Range("S18").Value = 0.1285714285714286
# Row 19
$ws.Range("F19").Value = 0.02178217821782178
$ws.Range("H19").Value = 0.2584158415841584
$ws.Range("I19").Value = 0.05346534653465346
$ws.Range("J19").Value = 0.3376237623762376
$ws.Range("K19").Value = 0.102970297029703
$ws.Range("M19").Value = 0.02574257425742574
$ws.Range("O19").Value = 0.07425742574257425
$ws.Range("S19").Value = 0.1257425742574257
